# Bang danh gia - cap nhat noi dung danh gia cho SV 1642076 (hang 8):
# thay "Cai dat DB" bang chi tiet cai dat Store Procedure, va dien cac
# o Dirty Read / Repeatable Read / Lost Update / Deadlock dang "Khong co"
# bang mo ta Victim/Cause cu the (kich ban Deadlock that).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ten chuc nang (cot E) - can doi canh trai nhu cot D/F
$ws.Range("E8").Value = "Database: Tranh chấp `nCài đặt Store Procedure:`n1. sp_GVQLCAPQUYEN`n2. sp_LayDanhSachGV"
$ws.Range("E8").HorizontalAlignment = -4131

# Dirty Read (GUI / SQL)
$ws.Range("G8").Value = "Giáo viên quản lý "
$ws.Range("H8").Value = "Victim: SP_LOADDSGV`nCause: SP_GVQLCAPQUYEN"

# Repeatable Read (GUI / SQL)
$ws.Range("I8").Value = "Giáo viên quản lý "
$ws.Range("J8").Value = "Victim: SP_LOAD_DSGV_CUNG_GVQL`nCause: SP_UPDATEGVQL"

# Lost Update (GUI / SQL)
$ws.Range("M8").Value = "Giáo viên quản lý "
$ws.Range("N8").Value = "Victim: sp_UpdateToanQuyenGV1`nCause: sp_UpdateToanQuyenGV2"

# Deadlock (GUI / SQL)
$ws.Range("O8").Value = "Giáo viên quản lý "
$ws.Range("P8").Value = "Victim: sp_UpdateToanQuyenGV2_Deadlock`nCause: sp_UpdateToanQuyenGV1_Deadlock"
